$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column W (23) width change: stored width 8 -> 9
$ws.Columns.Item(23).ColumnWidth = 8.17

# Cell value updates
$ws.Range("N15").Value = -0.864
$ws.Range("P15").Value = 0.714
$ws.Range("N16").Value = 0.095
$ws.Range("P16").Value = 11.305
$ws.Range("N17").Value = 0.095
$ws.Range("P17").Value = 0.555
$ws.Range("N18").Value = -38.405
$ws.Range("P18").Value = 48.705
$ws.Range("N19").Value = -27.401
$ws.Range("P19").Value = 27.401
$ws.Range("N20").Value = -0.864
$ws.Range("P20").Value = 0.864
$ws.Range("N21").Value = 0.095
$ws.Range("P21").Value = -0.095
$ws.Range("W21").Value = 0.095
$ws.Range("N22").Value = 0.095
$ws.Range("P22").Value = -0.095
$ws.Range("W22").Value = 0.095
$ws.Range("N25").Value = -1.631
$ws.Range("P25").Value = 2.031
$ws.Range("N26").Value = -0.864
$ws.Range("P26").Value = -1.286
$ws.Range("W29").Value = -39.364
$ws.Range("W30").Value = -28.364
$ws.Range("W31").Value = -28.364
$ws.Range("W32").Value = -39.364
$ws.Range("W33").Value = -39.364
$ws.Range("W34").Value = -39.364
$ws.Range("W35").Value = -28.364
$ws.Range("W36").Value = -28.364
$ws.Range("W37").Value = -17.364
$ws.Range("W38").Value = -17.364
$ws.Range("W39").Value = -6.364
$ws.Range("W40").Value = -6.364
$ws.Range("W41").Value = -6.364
$ws.Range("W42").Value = -17.364
$ws.Range("W43").Value = -17.364
$ws.Range("W44").Value = -6.364
$ws.Range("W45").Value = 4.636
$ws.Range("W46").Value = 4.636
$ws.Range("W47").Value = 15.636
$ws.Range("W48").Value = 15.636
$ws.Range("W49").Value = 26.636
$ws.Range("W50").Value = 37.636
$ws.Range("W51").Value = 37.636
$ws.Range("W52").Value = 26.636
$ws.Range("W53").Value = 26.636
$ws.Range("W54").Value = 37.636
$ws.Range("W55").Value = 37.636
$ws.Range("W56").Value = 26.636
$ws.Range("W57").Value = 15.636
$ws.Range("W58").Value = 15.636
$ws.Range("W59").Value = 4.636
$ws.Range("W60").Value = 4.636
$ws.Range("W61").Value = 4.636
$ws.Range("W62").Value = 4.636
$ws.Range("W63").Value = 15.636
$ws.Range("W64").Value = 15.636
$ws.Range("W65").Value = 26.636
$ws.Range("W66").Value = 37.636
$ws.Range("W67").Value = 37.636
$ws.Range("W68").Value = 26.636
$ws.Range("W69").Value = 26.636
$ws.Range("W70").Value = 37.636
$ws.Range("W71").Value = 37.636
$ws.Range("W72").Value = 26.636
$ws.Range("W73").Value = 15.636
$ws.Range("W74").Value = 15.636
$ws.Range("W75").Value = 4.636
$ws.Range("W76").Value = 4.636
$ws.Range("W77").Value = -6.364
$ws.Range("W78").Value = -17.364
$ws.Range("W79").Value = -17.364
$ws.Range("W80").Value = -6.364
$ws.Range("W81").Value = -6.364
$ws.Range("W82").Value = -6.364
$ws.Range("W83").Value = -17.364
$ws.Range("W84").Value = -17.364
$ws.Range("W85").Value = -28.364
$ws.Range("W86").Value = -28.364
$ws.Range("W87").Value = -39.364
$ws.Range("W88").Value = -39.364
$ws.Range("W89").Value = -39.364
$ws.Range("W90").Value = -28.364
$ws.Range("W91").Value = -28.364
$ws.Range("W92").Value = -39.364
$ws.Range("N93").Value = -0.864
$ws.Range("P93").Value = 0.864
$ws.Range("N94").Value = -0.864
$ws.Range("P94").Value = 1.264
$ws.Range("N95").Value = -0.864
$ws.Range("P95").Value = 0.864
$ws.Range("N98").Value = 0.095
$ws.Range("P98").Value = -0.445
$ws.Range("N99").Value = -0.108
$ws.Range("P99").Value = 0.758
$ws.Range("N100").Value = -1.526
$ws.Range("P100").Value = 1.276
$ws.Range("N101").Value = 0.095
$ws.Range("P101").Value = 11.905
$ws.Range("N102").Value = -38.405
$ws.Range("P102").Value = 49.205
$ws.Range("N103").Value = -27.401
$ws.Range("P103").Value = 28.201
$ws.Range("N104").Value = 0.095
$ws.Range("P104").Value = -0.095
$ws.Range("W104").Value = 0.095
$ws.Range("N105").Value = 0.095
$ws.Range("P105").Value = -0.095
$ws.Range("W105").Value = 0.095

Write-Host "Done. Applied 110 cell changes and 1 column width change."
